# Apply weekly refresh of Fruta / hortaliza data: rows 2-21 get their
# date / volume / price / origin figures re-permuted (values move between
# rows while row 13 stays put).  Only columns D, J, K, L, M, O, P change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row (source row's old values become the
# destination row's new values).
$mapping = @{
    2  = 17
    3  = 6
    4  = 9
    5  = 10
    6  = 11
    7  = 14
    8  = 4
    9  = 8
    10 = 3
    11 = 12
    12 = 18
    13 = 13
    14 = 15
    15 = 19
    16 = 20
    17 = 5
    18 = 16
    19 = 2
    20 = 21
    21 = 7
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot current values for the columns that change, before we start
# overwriting cells (since sources and destinations overlap).
$snapshot = @{}
for ($r = 2; $r -le 21; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $row
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
